$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 26.45917141421879
$ws.Range("C2").Value = 12.911434130392
$ws.Range("D2").Value = 4.391986222072624
$ws.Range("E2").Value = 9.749158777423458
$ws.Range("F2").Value = 50.56703541513724
$ws.Range("I2").Value = 34.96074594669077
$ws.Range("J2").Value = 9.434583285221256
$ws.Range("L2").Value = 12.84096026796663
$ws.Range("N2").Value = 20.87467228318587
$ws.Range("B3").Value = 26.13633795757265
$ws.Range("C3").Value = 12.5315896258594
$ws.Range("D3").Value = 4.373720773786265
$ws.Range("E3").Value = 9.757491529661477
$ws.Range("F3").Value = 50.52868156315316
$ws.Range("I3").Value = 35.01450159061653
$ws.Range("J3").Value = 9.452025722434126
$ws.Range("L3").Value = 12.84170914579781
$ws.Range("N3").Value = 20.94969735457105
$ws.Range("B4").Value = 25.9430940541769
$ws.Range("C4").Value = 12.29659075060337
$ws.Range("D4").Value = 4.362163248105515
$ws.Range("E4").Value = 9.762969755326878
$ws.Range("F4").Value = 50.51795759874397
$ws.Range("I4").Value = 35.05525172982829
$ws.Range("J4").Value = 9.463303736571183
$ws.Range("L4").Value = 12.84432891059916
$ws.Range("N4").Value = 20.99778645518981
$ws.Range("B5").Value = 25.86568777793827
$ws.Range("C5").Value = 12.20054596236009
$ws.Range("D5").Value = 4.357366874865952
$ws.Range("E5").Value = 9.765293369495872
$ws.Range("F5").Value = 50.51680918569031
$ws.Range("I5").Value = 35.07379794977471
$ws.Range("J5").Value = 9.468042985702576
$ws.Range("L5").Value = 12.84594017774399
$ws.Range("N5").Value = 21.0178935576011
$ws.Range("B6").Value = 25.85291818287793
$ws.Range("C6").Value = 12.18458582745287
$ws.Range("D6").Value = 4.356565197219108
$ws.Range("E6").Value = 9.765684717443893
$ws.Range("F6").Value = 50.51681292933466
$ws.Range("I6").Value = 35.07699450221219
$ws.Range("J6").Value = 9.468838607121603
$ws.Range("L6").Value = 12.84624057737366
$ws.Range("N6").Value = 21.02126320082472
$ws.Range("B7").Value = 25.94204457405184
$ws.Range("C7").Value = 12.29529636965986
$ws.Range("D7").Value = 4.362098913725484
$ws.Range("E7").Value = 9.763000722891348
$ws.Range("F7").Value = 50.51792907283051
$ws.Range("I7").Value = 35.05549400472915
$ws.Range("J7").Value = 9.463367070641075
$ws.Range("L7").Value = 12.84434843875106
$ws.Range("N7").Value = 20.99805555794628
$ws.Range("B8").Value = 26.34688061612846
$ws.Range("C8").Value = 12.78092785368473
$ws.Range("D8").Value = 4.385759088595742
$ws.Range("E8").Value = 9.751956918173788
$ws.Range("F8").Value = 50.5511477459254
$ws.Range("I8").Value = 34.97766952846975
$ws.Range("J8").Value = 9.440479770867212
$ws.Range("L8").Value = 12.84077057065531
$ws.Range("N8").Value = 20.90012194028181
$ws.Range("B9").Value = 27.17615630022667
$ws.Range("C9").Value = 13.7124166764653
$ws.Range("D9").Value = 4.429461493143267
$ws.Range("E9").Value = 9.733162822080379
$ws.Range("F9").Value = 50.71810602218246
$ws.Range("I9").Value = 34.88679777846943
$ws.Range("J9").Value = 9.400085652749766
$ws.Range("L9").Value = 12.85086184844289
$ws.Range("N9").Value = 20.72405446489158
$ws.Range("B10").Value = 27.80140797251795
$ws.Range("C10").Value = 14.37575874740175
$ws.Range("D10").Value = 4.459958681125066
$ws.Range("E10").Value = 9.721088383441003
$ws.Range("F10").Value = 50.9027692798671
$ws.Range("I10").Value = 34.85805098449528
$ws.Range("J10").Value = 9.373114147156398
$ws.Range("L10").Value = 12.8686542992509
$ws.Range("N10").Value = 20.60433261798775
$ws.Range("B11").Value = 28.08812335864507
$ws.Range("C11").Value = 14.67142718898914
$ws.Range("D11").Value = 4.473489748828265
$ws.Range("E11").Value = 9.715969427790396
$ws.Range("F11").Value = 51.00017328555111
$ws.Range("I11").Value = 34.85329714505836
$ws.Range("J11").Value = 9.361425346720827
$ws.Range("L11").Value = 12.87898710214178
$ws.Range("N11").Value = 20.55193731565147
$ws.Range("B12").Value = 28.19692419040041
$ws.Range("C12").Value = 14.78240060999072
$ws.Range("D12").Value = 4.478564908464574
$ws.Range("E12").Value = 9.714084577126371
$ws.Range("F12").Value = 51.03897447646604
$ws.Range("I12").Value = 34.85269835654319
$ws.Range("J12").Value = 9.357082118524785
$ws.Range("L12").Value = 12.88322013830551
$ws.Range("N12").Value = 20.53239209602318
$ws.Range("B13").Value = 28.17348343299061
$ws.Range("C13").Value = 14.75854622074563
$ws.Range("D13").Value = 4.477474043527012
$ws.Range("E13").Value = 9.714488132677655
$ws.Range("F13").Value = 51.0305329063627
$ws.Range("I13").Value = 34.85277381948777
$ws.Range("J13").Value = 9.358013822874463
$ws.Range("L13").Value = 12.88229426582725
$ws.Range("N13").Value = 20.53658837846378
$ws.Range("B14").Value = 28.09707047140556
$ws.Range("C14").Value = 14.68057749099425
$ws.Range("D14").Value = 4.473908259931664
$ws.Range("E14").Value = 9.715813286939442
$ws.Range("F14").Value = 51.00332714412004
$ws.Range("I14").Value = 34.85322378128846
$ws.Range("J14").Value = 9.36106636441156
$ws.Range("L14").Value = 12.87932894939851
$ws.Range("N14").Value = 20.55032339979951
$ws.Range("B15").Value = 28.05029204877444
$ws.Range("C15").Value = 14.63268729129263
$ws.Range("D15").Value = 4.471717777719905
$ws.Range("E15").Value = 9.716631956063551
$ws.Range("F15").Value = 50.98691204843424
$ws.Range("I15").Value = 34.85365597461548
$ws.Range("J15").Value = 9.362946939188927
$ws.Range("L15").Value = 12.87755425560871
$ws.Range("N15").Value = 20.55877496865302
$ws.Range("B16").Value = 27.78270779303389
$ws.Range("C16").Value = 14.35630357535206
$ws.Range("D16").Value = 4.459067545798549
$ws.Range("E16").Value = 9.7214304265792
$ws.Range("F16").Value = 50.89667250003318
$ws.Range("I16").Value = 34.85852954502906
$ws.Range("J16").Value = 9.373889685281387
$ws.Range("L16").Value = 12.86802393658911
$ws.Range("N16").Value = 20.60779823820624
$ws.Range("B17").Value = 27.61906675171904
$ws.Range("C17").Value = 14.18510675610471
$ws.Range("D17").Value = 4.45121970171731
$ws.Range("E17").Value = 9.724469748199034
$ws.Range("F17").Value = 50.8447388223658
$ws.Range("I17").Value = 34.86365452496342
$ws.Range("J17").Value = 9.380751120683795
$ws.Range("L17").Value = 12.86274967831322
$ws.Range("N17").Value = 20.63840073544587
$ws.Range("B18").Value = 27.5251658055536
$ws.Range("C18").Value = 14.08607184114632
$ws.Range("D18").Value = 4.446673587546247
$ws.Range("E18").Value = 9.726253075710142
$ws.Range("F18").Value = 50.81612998235821
$ws.Range("I18").Value = 34.86738547070287
$ws.Range("J18").Value = 9.38475231825438
$ws.Range("L18").Value = 12.859926920255
$ws.Range("N18").Value = 20.65619710242625
$ws.Range("B19").Value = 27.49341361576426
$ws.Range("C19").Value = 14.0524465934396
$ws.Range("D19").Value = 4.445128799998569
$ws.Range("E19").Value = 9.726862928648007
$ws.Range("F19").Value = 50.80666055396046
$ws.Range("I19").Value = 34.86878307023759
$ws.Range("J19").Value = 9.386116459822686
$ws.Range("L19").Value = 12.85900744658411
$ws.Range("N19").Value = 20.6622561171595
$ws.Range("B20").Value = 27.63646447589992
$ws.Range("C20").Value = 14.20339049778655
$ws.Range("D20").Value = 4.452058451728844
$ws.Range("E20").Value = 9.724142566356461
$ws.Range("F20").Value = 50.85013670267128
$ws.Range("I20").Value = 34.86302786979215
$ws.Range("J20").Value = 9.380015053002429
$ws.Range("L20").Value = 12.86328932115689
$ws.Range("N20").Value = 20.63512291553068
$ws.Range("B21").Value = 28.11950941885351
$ws.Range("C21").Value = 14.70350651508585
$ws.Range("D21").Value = 4.474956936786468
$ws.Range("E21").Value = 9.715422603943761
$ws.Range("F21").Value = 51.01126622108506
$ws.Range("I21").Value = 34.8530589785993
$ws.Range("J21").Value = 9.360167507995619
$ws.Range("L21").Value = 12.88019125866607
$ws.Range("N21").Value = 20.54628107619444
$ws.Range("B22").Value = 28.43649034929707
$ws.Range("C22").Value = 15.02454666735179
$ws.Range("D22").Value = 4.489638407942845
$ws.Range("E22").Value = 9.710035877064289
$ws.Range("F22").Value = 51.1277373424344
$ws.Range("I22").Value = 34.85354795610859
$ws.Range("J22").Value = 9.347679968060364
$ws.Range("L22").Value = 12.89310323051206
$ws.Range("N22").Value = 20.48994095713605
$ws.Range("B23").Value = 28.26722707845247
$ws.Range("C23").Value = 14.85376841875851
$ws.Range("D23").Value = 4.481828467962293
$ws.Range("E23").Value = 9.712882353708835
$ws.Range("F23").Value = 51.06455714276813
$ws.Range("I23").Value = 34.85264477363719
$ws.Range("J23").Value = 9.354300661518678
$ws.Range("L23").Value = 12.88604179447055
$ws.Range("N23").Value = 20.51985354667938
$ws.Range("B24").Value = 27.62859840125601
$ws.Range("C24").Value = 14.19512631433404
$ws.Range("D24").Value = 4.451679359670043
$ws.Range("E24").Value = 9.724290373182383
$ws.Range("F24").Value = 50.84769243135162
$ws.Range("I24").Value = 34.86330873732054
$ws.Range("J24").Value = 9.380347653286503
$ws.Range("L24").Value = 12.8630446960054
$ws.Range("N24").Value = 20.63660418680136
$ws.Range("B25").Value = 26.94860057769209
$ws.Range("C25").Value = 13.46355644556022
$ws.Range("D25").Value = 4.417924044896367
$ws.Range("E25").Value = 9.73794185939531
$ws.Range("F25").Value = 50.66203660410969
$ws.Range("I25").Value = 34.9047298226711
$ws.Range("J25").Value = 9.410535975687189
$ws.Range("L25").Value = 12.84630481507732
$ws.Range("N25").Value = 20.76998570992138
